$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1
$ws.Range("AD2").Value = 0
$ws.Range("AH2").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 4
$ws.Range("AD5").Value = 4
$ws.Range("AE5").ClearContents()
$ws.Range("X6").Value = 4
$ws.Range("X7").Value = 0
$ws.Range("AD16").Value = 4
$ws.Range("AE16").Value = 0
$ws.Range("Z17").Value = 0
$ws.Range("AA17").Value = 4
$ws.Range("AD19").Value = 0
$ws.Range("AH19").Value = 1

$ws.Range("AJ18").Select()
